$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column H ("Is verified") ---------------------------------------
# Copy the formatting (fill/border/font/number-format) from column G onto
# column H for every row so the new column matches the existing style
# pattern (header/data rows styled, spacer rows styled differently, etc.)
$ws.Range("G1:G50").Copy()
$ws.Range("H1:H50").PasteSpecial(-4122)

# Match column F/G's width for the new column H.
$ws.Range("H1").EntireColumn.ColumnWidth = $ws.Range("F1").EntireColumn.ColumnWidth

# Header
$ws.Range("H1").Value = "Is verified"

# Data rows
$ws.Range("H2").Value = "Yes"
$ws.Range("H3").Value = "No"
$ws.Range("H5").Value = "yes"
$ws.Range("H6").Value = "no"

# --- Fix the typo'd e-mail address in row 5 -----------------------------
$ws.Range("F5").Value = "ttest@example.com"

# Rebuild the hyperlinks so F5's display text matches the corrected
# address (the COM layer can only append new hyperlink entries, so the
# whole collection is cleared and re-created in order).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:test@example.com", "", "", "test@example.com")
$ws.Hyperlinks.Add($ws.Range("F3"), "mailto:test@example.com", "", "", "test1@example.com")
$ws.Hyperlinks.Add($ws.Range("F4"), "mailto:test@example.com", "", "", "test2@example.com")
$ws.Hyperlinks.Add($ws.Range("F5"), "mailto:test@example.com", "", "", "ttest@example.com")
$ws.Hyperlinks.Add($ws.Range("F6"), "mailto:test@example.com", "", "", "test3@example.com")
